$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 5.277653829439964
$ws.Range("D2").Value = 6.035145143649665
$ws.Range("E2").Value = 16.37091983698481
$ws.Range("F2").Value = 34.68142715939951
$ws.Range("G2").Value = 3.655314214704473
$ws.Range("I2").Value = 20.78047897258175
$ws.Range("K2").Value = 16.83311981035276
$ws.Range("N2").Value = 18.98595218757068
$ws.Range("B3").Value = 5.211296798088775
$ws.Range("D3").Value = 6.062964248149949
$ws.Range("E3").Value = 15.44472284728718
$ws.Range("F3").Value = 34.27243119398937
$ws.Range("G3").Value = 3.659784321781402
$ws.Range("I3").Value = 20.83734276693864
$ws.Range("K3").Value = 16.30011171034324
$ws.Range("N3").Value = 19.04901790168128
$ws.Range("B4").Value = 5.171183486577794
$ws.Range("D4").Value = 6.08108398883825
$ws.Range("E4").Value = 14.85290666453555
$ws.Range("F4").Value = 34.03120406656212
$ws.Range("G4").Value = 3.662664203040523
$ws.Range("I4").Value = 20.87566063376613
$ws.Range("K4").Value = 15.96980556576823
$ws.Range("N4").Value = 19.08965856504325
$ws.Range("B5").Value = 5.155016074071039
$ws.Range("D5").Value = 6.088728162740775
$ws.Range("E5").Value = 14.6061952229365
$ws.Range("F5").Value = 33.93549150400312
$ws.Range("G5").Value = 3.663871938778475
$ws.Range("I5").Value = 20.89212869433245
$ws.Range("K5").Value = 15.83467296391784
$ws.Range("N5").Value = 19.10670284144036
$ws.Range("B6").Value = 5.152342902761855
$ws.Range("D6").Value = 6.090013162375148
$ws.Range("E6").Value = 14.56490333879608
$ws.Range("F6").Value = 33.91975770802576
$ws.Range("G6").Value = 3.664074550261617
$ws.Range("I6").Value = 20.89491466049785
$ws.Range("K6").Value = 15.81220938232984
$ws.Range("N6").Value = 19.10956221202376
$ws.Range("B7").Value = 5.170964694971278
$ws.Range("D7").Value = 6.081186028384698
$ws.Range("E7").Value = 14.84960146422651
$ws.Range("F7").Value = 34.02990264093963
$ws.Range("G7").Value = 3.662680352459496
$ws.Range("I7").Value = 20.87587927649637
$ws.Range("K7").Value = 15.96798494767828
$ws.Range("N7").Value = 19.0898864741819
$ws.Range("B8").Value = 5.254655558426209
$ws.Range("D8").Value = 6.044520999859149
$ws.Range("E8").Value = 16.05651588964499
$ws.Range("F8").Value = 34.53841327475705
$ws.Range("G8").Value = 3.656827539420501
$ws.Range("I8").Value = 20.7993777230609
$ws.Range("K8").Value = 16.65010447120298
$ws.Range("N8").Value = 19.00729949110637
$ws.Range("B9").Value = 5.422861207195067
$ws.Range("D9").Value = 5.980907331573978
$ws.Range("E9").Value = 18.27871368266245
$ws.Range("F9").Value = 35.60942583136313
$ws.Range("G9").Value = 3.646415794325011
$ws.Range("I9").Value = 20.67647666030679
$ws.Range("K9").Value = 17.95416078343207
$ws.Range("N9").Value = 18.86053905564573
$ws.Range("B10").Value = 5.547714917893387
$ws.Range("D10").Value = 5.939286451511165
$ws.Range("E10").Value = 19.92421650665488
$ws.Range("F10").Value = 36.43466118560521
$ws.Range("G10").Value = 3.639405671770794
$ws.Range("I10").Value = 20.60287113126041
$ws.Range("K10").Value = 18.88025202998686
$ws.Range("N10").Value = 18.76193766984286
$ws.Range("B11").Value = 5.604536199992551
$ws.Range("D11").Value = 5.921478209014609
$ws.Range("E11").Value = 20.63177144703029
$ws.Range("F11").Value = 36.81693329932443
$ws.Range("G11").Value = 3.63635321692425
$ws.Range("I11").Value = 20.57304402476769
$ws.Range("K11").Value = 19.29251056768478
$ws.Range("N11").Value = 18.71907652431376
$ws.Range("B12").Value = 5.626035808526846
$ws.Range("D12").Value = 5.91489792362874
$ws.Range("E12").Value = 20.89385099683383
$ws.Range("F12").Value = 36.96255136019669
$ws.Range("G12").Value = 3.635216785896709
$ws.Range("I12").Value = 20.56227786776939
$ws.Range("K12").Value = 19.44717255575911
$ws.Range("N12").Value = 18.7031322969161
$ws.Range("B13").Value = 5.621406583137868
$ws.Range("D13").Value = 5.916307817922292
$ws.Range("E13").Value = 20.83766731816026
$ws.Range("F13").Value = 36.93115369399797
$ws.Range("G13").Value = 3.635460673370595
$ws.Range("I13").Value = 20.56457299065403
$ws.Range("K13").Value = 19.41393015224205
$ws.Range("N13").Value = 18.70655343923013
$ws.Range("B14").Value = 5.60630545183751
$ws.Range("D14").Value = 5.920933563175758
$ws.Range("E14").Value = 20.65345003489412
$ws.Range("F14").Value = 36.82889688265484
$ws.Range("G14").Value = 3.636259332817614
$ws.Range("I14").Value = 20.5721476734384
$ws.Range("K14").Value = 19.30526461206707
$ws.Range("N14").Value = 18.71775904538726
$ws.Range("B15").Value = 5.597052698373142
$ws.Range("D15").Value = 5.923788279103307
$ws.Range("E15").Value = 20.53985012807187
$ws.Range("F15").Value = 36.76636982666655
$ws.Range("G15").Value = 3.636751065288478
$ws.Range("I15").Value = 20.57685632729921
$ws.Range("K15").Value = 19.23851053956145
$ws.Range("N15").Value = 18.72466008822704
$ws.Range("B16").Value = 5.544000208650301
$ws.Range("D16").Value = 5.940473031729391
$ws.Range("E16").Value = 19.87715454824565
$ws.Range("F16").Value = 36.40980609587759
$ws.Range("G16").Value = 3.639607889718899
$ws.Range("I16").Value = 20.60489423013096
$ws.Range("K16").Value = 18.85311485525395
$ws.Range("N16").Value = 18.76477883090288
$ws.Range("B17").Value = 5.511445486175597
$ws.Range("D17").Value = 5.950997840419913
$ws.Range("E17").Value = 19.46014054666395
$ws.Range("F17").Value = 36.1927341446485
$ws.Range("G17").Value = 3.6413953030597
$ws.Range("I17").Value = 20.62303316037535
$ws.Range("K17").Value = 18.61426105799048
$ws.Range("N17").Value = 18.78990067558883
$ws.Range("B18").Value = 5.492724416146457
$ws.Range("D18").Value = 5.957157256835075
$ws.Range("E18").Value = 19.21642441205133
$ws.Range("F18").Value = 36.06853428044942
$ws.Range("G18").Value = 3.642436232092302
$ws.Range("I18").Value = 20.63381017787897
$ws.Range("K18").Value = 18.47603731097127
$ws.Range("N18").Value = 18.8045377106797
$ws.Range("B19").Value = 5.486387039209884
$ws.Range("D19").Value = 5.959260860733815
$ws.Range("E19").Value = 19.13324206084259
$ws.Range("F19").Value = 36.02659858583605
$ws.Range("G19").Value = 3.642790885669204
$ws.Range("I19").Value = 20.63751808095156
$ws.Range("K19").Value = 18.42909756984064
$ws.Range("N19").Value = 18.80952579040297
$ws.Range("B20").Value = 5.51491077843472
$ws.Range("D20").Value = 5.949866493557508
$ws.Range("E20").Value = 19.50493168836261
$ws.Range("F20").Value = 36.21577504439822
$ws.Range("G20").Value = 3.641203700445945
$ws.Range("I20").Value = 20.62106661793608
$ws.Range("K20").Value = 18.63977569459344
$ws.Range("N20").Value = 18.78720699689461
$ws.Range("B21").Value = 5.610741655731854
$ws.Range("D21").Value = 5.919570425185513
$ws.Range("E21").Value = 20.70771774418387
$ws.Range("F21").Value = 36.85890986717671
$ws.Range("G21").Value = 3.636024220107652
$ws.Range("I21").Value = 20.56990843281834
$ws.Range("K21").Value = 19.33722284845427
$ws.Range("N21").Value = 18.7144599173417
$ws.Range("B22").Value = 5.673262100030112
$ws.Range("D22").Value = 5.900722837519683
$ws.Range("E22").Value = 21.45969165214097
$ws.Range("F22").Value = 37.28418884324787
$ws.Range("G22").Value = 3.632752526940036
$ws.Range("I22").Value = 20.53955674612192
$ws.Range("K22").Value = 19.78452525284106
$ws.Range("N22").Value = 18.66858452203111
$ws.Range("B23").Value = 5.639910641204317
$ws.Range("D23").Value = 5.910694468232988
$ws.Range("E23").Value = 21.06145796536646
$ws.Range("F23").Value = 37.05679837967813
$ws.Range("G23").Value = 3.634488368005819
$ws.Range("I23").Value = 20.55547293772428
$ws.Range("K23").Value = 19.54661765741231
$ws.Range("N23").Value = 18.69291644017711
$ws.Range("B24").Value = 5.513344133555376
$ws.Range("D24").Value = 5.950377637434237
$ws.Range("E24").Value = 19.48469397938399
$ws.Range("F24").Value = 36.20535638225326
$ws.Range("G24").Value = 3.641290282482354
$ws.Range("I24").Value = 20.62195460585708
$ws.Range("K24").Value = 18.62824332981029
$ws.Range("N24").Value = 18.78842420427308
$ws.Range("B25").Value = 5.377046013934224
$ws.Range("D25").Value = 5.997223004385758
$ws.Range("E25").Value = 17.66398126275979
$ws.Range("F25").Value = 35.31245620135013
$ws.Range("G25").Value = 3.649119423963286
$ws.Range("I25").Value = 20.70680513913403
$ws.Range("K25").Value = 17.6062299993847
$ws.Range("N25").Value = 18.89861926794912
